$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.933.06"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "1.630.20"

$ws.Range("E4").Value = "  +0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("E7").Value = "  +0.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.28"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.79%  "

$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0907"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").Value = "1.863.92"

$ws.Range("D13").Value = "1.636.37"
$ws.Range("E13").Value = "  +2.11%  "

$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "29.933.95"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +13.45%  "

$ws.Range("E17").Value = "  +0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.45"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("E22").Value = "  +1.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.73"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.50"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("E28").Value = "  +1.79%  "

$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.22%  "

$ws.Range("E32").Value = "  +3.54%  "

$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("D34").Value = "1.423.73"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("E35").Value = "  +4.16%  "

$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.73"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.06%  "

$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.26"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.549"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("E42").Value = "  +1.61%  "

$ws.Range("E43").Value = "  +0.52%  "

$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.02"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.58"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.70%  "

$ws.Range("D49").Value = "1.771.63"
$ws.Range("E49").Value = "  +1.71%  "

$ws.Range("E50").Value = "  +10.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.29"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.98%  "
